$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Multiply the data block B2:J16 by 10 (values were stored as fractions
# one decade too small, e.g. 1.9E-3 -> 1.9E-2)
$rng = $ws.Range("B2:J16")
foreach ($cell in $rng) {
    $cell.Value = $cell.Value() * 10
}

# Update the view: show zoom-normal percentage and move the active selection
$excel.ActiveWindow.Zoom = 100
[void]$ws.Range("K24").Select()
